$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.142.30'
$ws.Range("E2").Value = '  +0.29%  '
$ws.Range("D3").Value = '1.875.43'
$ws.Range("E3").Value = '  -0.90%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9992'
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '243.00'
$ws.Range("E5").Value = '  -1.56%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9990'
$ws.Range("E6").Value = '  +0.15%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4902'
$ws.Range("E7").Value = '  -1.49%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2911'
$ws.Range("E8").Value = '  -1.33%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06591'
$ws.Range("E9").Value = '  -0.82%  '
$ws.Range("D10").Value = '1.877.47'
$ws.Range("E10").Value = '  -0.78%  '
$ws.Range("E11").Value = '  -3.47%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07192'
$ws.Range("E12").Value = '  -0.40%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.6662'
$ws.Range("E13").Value = '  -1.52%  '
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.918'
$ws.Range("E14").Value = '  +1.61%  '
$ws.Range("B15").Value = 'Litecoin'
$ws.Range("C15").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '86.01'
$ws.Range("E15").Value = '  +0.21%  '
$ws.Range("D16").Value = '30.097.10'
$ws.Range("E16").Value = '  +0.16%  '
$ws.Range("B17").Value = 'ShibaInu'
$ws.Range("C17").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000007789'
$ws.Range("E17").Value = '  -2.95%  '
$ws.Range("B18").Value = 'Dai'
$ws.Range("C18").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.9993'
$ws.Range("E18").Value = '  +0.14%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.77'
$ws.Range("E19").Value = '  -1.03%  '
$ws.Range("E20").Value = '  -0.48%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9981'
$ws.Range("E21").Value = '  +0.03%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.768'
$ws.Range("E22").Value = '  -0.19%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.822'
$ws.Range("E23").Value = '  +2.73%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.170'
$ws.Range("E24").Value = '  +0.02%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '152.79'
$ws.Range("E25").Value = '  +3.90%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '143.63'
$ws.Range("E26").Value = '  +8.33%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '16.96'
$ws.Range("E27").Value = '  +0.71%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.890'
$ws.Range("E28").Value = '  -3.35%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.396'
$ws.Range("E29").Value = '  +2.36%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.213'
$ws.Range("E30").Value = '  -0.62%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08790'
$ws.Range("E31").Value = '  +0.55%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.989'
$ws.Range("E32").Value = '  +0.91%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05140'
$ws.Range("E33").Value = '  +0.71%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7162'
$ws.Range("E34").Value = '  +1.49%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.108'
$ws.Range("E35").Value = '  -0.91%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.664'
$ws.Range("E36").Value = '  +0.11%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.01843'
$ws.Range("E37").Value = '  +11.03%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.680'
$ws.Range("E38").Value = '  -4.05%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.141'
$ws.Range("E39").Value = '  -3.79%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.9296'
$ws.Range("E40").Value = '  -1.31%  '
$ws.Range("B41").Value = 'FraxShare'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.782'
$ws.Range("E41").Value = '  -4.95%  '
$ws.Range("B42").Value = 'PaxDollar'
$ws.Range("C42").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9987'
$ws.Range("E42").Value = '  +0.27%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.4224'
$ws.Range("E43").Value = '  +0.11%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '103.50'
$ws.Range("E44").Value = '  +0.29%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '7.383'
$ws.Range("E45").Value = '  -1.48%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.1282'
$ws.Range("E46").Value = '  +1.68%  '
$ws.Range("E47").Value = '  -0.63%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '32.81'
$ws.Range("E48").Value = '  -0.21%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.3755'
$ws.Range("E49").Value = '  +0.29%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.234'
$ws.Range("E50").Value = '  -0.84%  '
$ws.Range("B51").Value = 'NEARProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.337'
$ws.Range("E51").Value = '  -0.42%  '
